$wb = $excel.ActiveWorkbook

# --- Sheet: LP1912 ---
$ws = $wb.Worksheets.Item("LP1912")
$ws.Cells.Item(2,1).Value = "Última actualización: 06:35:33"
$ws.Cells.Item(3,1).Value = "Total filas: 66"
$ws.Cells.Item(8,1).Value = "03:52:04"
$ws.Cells.Item(8,3).Value = "215A_EL PATO"
$ws.Cells.Item(8,4).Value = 54
$ws.Cells.Item(9,3).Value = "15_ABASTO"
$ws.Cells.Item(10,1).Value = "04:44:46"
$ws.Cells.Item(10,3).Value = "215_EL PELIGRO"
$ws.Cells.Item(10,4).Value = 2
$ws.Cells.Item(38,3).Value = "17X38_ROMERO"
$ws.Cells.Item(39,3).Value = "16_SANTA ANA"
$ws.Cells.Item(41,1).Value = "06:35:33"
$ws.Cells.Item(41,4).Value = 4
$ws.Cells.Item(42,1).Value = "06:35:33"
$ws.Cells.Item(42,2).Value = "06:41"
$ws.Cells.Item(42,4).Value = 6
$ws.Cells.Item(43,1).Value = "06:18:01"
$ws.Cells.Item(43,2).Value = "06:45"
$ws.Cells.Item(43,4).Value = 27
$ws.Cells.Item(44,1).Value = "05:16:02"
$ws.Cells.Item(44,2).Value = "06:50"
$ws.Cells.Item(44,3).Value = "17_ROMERO"
$ws.Cells.Item(44,4).Value = 94
$ws.Cells.Item(45,1).Value = "06:35:33"
$ws.Cells.Item(45,2).Value = "06:51"
$ws.Cells.Item(45,3).Value = "215A_EL PATO"
$ws.Cells.Item(45,4).Value = 16
$ws.Cells.Item(46,1).Value = "06:35:33"
$ws.Cells.Item(46,2).Value = "06:54"
$ws.Cells.Item(46,3).Value = "14_ABASTO"
$ws.Cells.Item(46,4).Value = 19
$ws.Cells.Item(47,1).Value = "06:35:33"
$ws.Cells.Item(47,2).Value = "07:04"
$ws.Cells.Item(47,3).Value = "225_GOMEZ"
$ws.Cells.Item(47,4).Value = 29
$ws.Cells.Item(48,1).Value = "06:35:33"
$ws.Cells.Item(48,2).Value = "07:06"
$ws.Cells.Item(48,4).Value = 31
$ws.Cells.Item(49,1).Value = "06:18:01"
$ws.Cells.Item(49,2).Value = "07:07"
$ws.Cells.Item(49,3).Value = "215C_EL PATO"
$ws.Cells.Item(49,4).Value = 49
$ws.Cells.Item(50,1).Value = "06:35:33"
$ws.Cells.Item(50,2).Value = "07:13"
$ws.Cells.Item(50,4).Value = 38
$ws.Cells.Item(51,1).Value = "06:18:01"
$ws.Cells.Item(51,2).Value = "07:14"
$ws.Cells.Item(51,3).Value = "14X44_ABASTO"
$ws.Cells.Item(51,4).Value = 56
$ws.Cells.Item(52,1).Value = "05:57:38"
$ws.Cells.Item(52,2).Value = "07:20"
$ws.Cells.Item(52,4).Value = 83
$ws.Cells.Item(53,1).Value = "06:35:33"
$ws.Cells.Item(53,2).Value = "07:21"
$ws.Cells.Item(53,3).Value = "215A_EL PATO"
$ws.Cells.Item(53,4).Value = 46
$ws.Cells.Item(54,1).Value = "06:35:33"
$ws.Cells.Item(54,2).Value = "07:24"
$ws.Cells.Item(54,3).Value = "16_SANTA ANA"
$ws.Cells.Item(54,4).Value = 49
$ws.Cells.Item(55,1).Value = "06:35:33"
$ws.Cells.Item(55,2).Value = "07:29"
$ws.Cells.Item(55,3).Value = "14_ABASTO"
$ws.Cells.Item(55,4).Value = 54
$ws.Cells.Item(56,1).Value = "06:35:33"
$ws.Cells.Item(56,2).Value = "07:33"
$ws.Cells.Item(56,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(56,4).Value = 58
$ws.Cells.Item(57,1).Value = "06:35:33"
$ws.Cells.Item(57,2).Value = "07:36"
$ws.Cells.Item(57,4).Value = 61
$ws.Cells.Item(58,1).Value = "06:35:33"
$ws.Cells.Item(58,2).Value = "07:36"
$ws.Cells.Item(58,3).Value = "17X38_ROMERO"
$ws.Cells.Item(58,4).Value = 61
$ws.Cells.Item(59,2).Value = "07:37"
$ws.Cells.Item(59,3).Value = "27_EL RETIRO"
$ws.Cells.Item(59,4).Value = 79
$ws.Cells.Item(60,1).Value = "06:35:33"
$ws.Cells.Item(60,2).Value = "07:43"
$ws.Cells.Item(60,3).Value = "10_OLMOS"
$ws.Cells.Item(60,4).Value = 68
$ws.Cells.Item(61,2).Value = "07:44"
$ws.Cells.Item(61,3).Value = "10_OLMOS"
$ws.Cells.Item(61,4).Value = 86
$ws.Cells.Item(62,1).Value = "06:35:33"
$ws.Cells.Item(62,2).Value = "07:49"
$ws.Cells.Item(62,3).Value = "15_ABASTO"
$ws.Cells.Item(62,4).Value = 74
$ws.Cells.Item(63,1).Value = "06:35:33"
$ws.Cells.Item(63,2).Value = "07:58"
$ws.Cells.Item(63,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(63,4).Value = 83
$ws.Cells.Item(64,2).Value = "07:59"
$ws.Cells.Item(64,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(64,4).Value = 101
$ws.Cells.Item(65,1).Value = "06:35:33"
$ws.Cells.Item(65,2).Value = "07:59"
$ws.Cells.Item(65,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(65,4).Value = 84
$ws.Cells.Item(65,5).Value = "LP1912"
$ws.Cells.Item(66,1).Value = "06:18:01"
$ws.Cells.Item(66,2).Value = "08:00"
$ws.Cells.Item(66,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(66,4).Value = 102
$ws.Cells.Item(66,5).Value = "LP1912"
$ws.Cells.Item(67,1).Value = "06:35:33"
$ws.Cells.Item(67,2).Value = "08:01"
$ws.Cells.Item(67,3).Value = "16_SANTA ANA"
$ws.Cells.Item(67,4).Value = 86
$ws.Cells.Item(67,5).Value = "LP1912"
$ws.Cells.Item(68,1).Value = "06:35:33"
$ws.Cells.Item(68,2).Value = "08:03"
$ws.Cells.Item(68,3).Value = "17X38_ROMERO"
$ws.Cells.Item(68,4).Value = 88
$ws.Cells.Item(68,5).Value = "LP1912"
$ws.Cells.Item(69,1).Value = "06:35:33"
$ws.Cells.Item(69,2).Value = "08:14"
$ws.Cells.Item(69,3).Value = "10_OLMOS"
$ws.Cells.Item(69,4).Value = 99
$ws.Cells.Item(69,5).Value = "LP1912"
$ws.Cells.Item(70,1).Value = "06:35:33"
$ws.Cells.Item(70,2).Value = "08:19"
$ws.Cells.Item(70,3).Value = "17_ROMERO"
$ws.Cells.Item(70,4).Value = 104
$ws.Cells.Item(70,5).Value = "LP1912"
$ws.Cells.Item(71,1).Value = "06:35:33"
$ws.Cells.Item(71,2).Value = "08:34"
$ws.Cells.Item(71,3).Value = "215C_EL PATO"
$ws.Cells.Item(71,4).Value = 119
$ws.Cells.Item(71,5).Value = "LP1912"

# --- Sheet: LP1912-215 ---
$ws = $wb.Worksheets.Item("LP1912-215")
$ws.Cells.Item(2,1).Value = "Última actualización: 06:35:33"
$ws.Cells.Item(3,1).Value = "Total filas: 12"
$ws.Cells.Item(12,1).Value = "06:35:33"
$ws.Cells.Item(12,4).Value = 16
$ws.Cells.Item(13,1).Value = "06:35:33"
$ws.Cells.Item(13,4).Value = 31
$ws.Cells.Item(16,1).Value = "06:35:33"
$ws.Cells.Item(16,4).Value = 46
$ws.Cells.Item(17,1).Value = "06:35:33"
$ws.Cells.Item(17,2).Value = "08:34"
$ws.Cells.Item(17,3).Value = "215C_EL PATO"
$ws.Cells.Item(17,4).Value = 119
$ws.Cells.Item(17,5).Value = "LP1912"

# --- Sheet: 6203-6173 ---
$ws = $wb.Worksheets.Item("6203-6173")
$ws.Cells.Item(2,1).Value = "Última actualización: 06:35:33"
$ws.Cells.Item(3,1).Value = "Total filas: 4"
$ws.Cells.Item(7,1).Value = "06:35:33"
$ws.Cells.Item(7,4).Value = 52
$ws.Cells.Item(8,1).Value = "06:35:33"
$ws.Cells.Item(8,4).Value = 95
$ws.Cells.Item(9,1).Value = "06:35:33"
$ws.Cells.Item(9,2).Value = "08:23"
$ws.Cells.Item(9,3).Value = "215C_LA PLATA"
$ws.Cells.Item(9,4).Value = 108
$ws.Cells.Item(9,5).Value = "L6203"
